# Update "想去人数" (want-to-go count, column F) values on both the
# "展览" and "全部类型" sheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

$sheetExpo = $wb.Worksheets.Item("展览")
$sheetExpo.Range("F2").Value = 4296
$sheetExpo.Range("F3").Value = 2435
$sheetExpo.Range("F5").Value = 22
$sheetExpo.Range("F9").Value = 127
$sheetExpo.Range("F10").Value = 136
$sheetExpo.Range("F12").Value = 1598
$sheetExpo.Range("F14").Value = 3340

$sheetAll = $wb.Worksheets.Item("全部类型")
$sheetAll.Range("F2").Value = 4296
$sheetAll.Range("F3").Value = 2435
$sheetAll.Range("F5").Value = 22
$sheetAll.Range("F11").Value = 127
$sheetAll.Range("F12").Value = 136
$sheetAll.Range("F16").Value = 1598
$sheetAll.Range("F18").Value = 3340
